# "Mise en commun des AutoEvals dans un seul fichier"
#
# The "students" sheet used to be the active tab; a new "Visa Enseignant"
# config row (value "GGZ") is inserted into the "configs" table right
# after the "Enseignant" row, and the "configs" sheet becomes the active
# tab with cell A4 selected.

$wb = $excel.ActiveWorkbook
$wsConfigs = $wb.Worksheets.Item("configs")

# Insert a new blank row at row 4 (pushes Date debut..Evaluation numero
# down by one row, carrying their styles along).
$wsConfigs.Rows("4:4").Insert()

# Fill in the new row with the "Visa Enseignant" / "GGZ" pair.
$wsConfigs.Range("A4").Value = "Visa Enseignant"
$wsConfigs.Range("B4").Value = "GGZ"

# Grow the structured table ("Tableau1") so it covers the new row too.
$tbl = $wsConfigs.ListObjects.Item("Tableau1")
$tbl.Resize($wsConfigs.Range("A1:B10"))

# "configs" becomes the active sheet, with A4 selected.
$wsConfigs.Activate()
$wsConfigs.Range("A4").Select()
